# Updates to responsibility parser based on PR feedback
#
# Row 135's responsibility text (column G) was truncated - extend it to the
# full sentence, and add four new responsibility rows (136-139) that were
# parsed out of the same DoDI 6440.02 CH 1.pdf document/paragraph.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Smart punctuation used in some of the cell text below.
$rsquo  = [char]0x2019   # ’
$ldquo  = [char]0x201C   # “
$rdquo  = [char]0x201D   # ”

# --- Fix up the truncated text in row 135 -----------------------------
$ws.Cells.Item(135, 7).Value = "Implement CLIP requirements within their respective Department" + $rsquo + "s Active and Reserve Components and facilities under their supervision to include oversight, inspections, proficiency testing, personnel standards, and training in laboratories performing testing on human specimens as defined under " + $ldquo + "laboratory" + $rdquo + " in the Glossary of this instruction."

# --- Common values shared by the new rows ------------------------------
$filename  = "DoDI 6440.02 CH 1.pdf"
$docTitle  = "Clinical Laboratory Improvement Program (CLIP)"
# Leading apostrophe forces this to be stored as text (matching the rest of
# the column) instead of being auto-converted to the number 2.
$orgNum    = "'2."
$orgText   = "SECRETARIES OF THE MILITARY DEPARTMENTS.  The Secretaries of the Military Departments:"
$orgEnt    = "Military Departments"

# --- Row 136 -------------------------------------------------------------
$ws.Cells.Item(136, 1).Value = $filename
$ws.Cells.Item(136, 2).Value = $docTitle
$ws.Cells.Item(136, 3).Value = $orgNum
$ws.Cells.Item(136, 4).Value = $orgText
$ws.Cells.Item(136, 5).Value = $orgEnt
$ws.Cells.Item(136, 6).Value = "b."
$ws.Cells.Item(136, 7).Value = "In accordance with DoD Manual 6440.02 (Reference (g)), follow CLIP procedures for corrective action on laboratory facilities whose proficiency testing or performance criteria fall outside the standards of CLIP policy."
$ws.Cells.Item(136, 8).Value = "DoD"
$ws.Cells.Item(136, 9).Value = "active"

# --- Row 137 -------------------------------------------------------------
$ws.Cells.Item(137, 1).Value = $filename
$ws.Cells.Item(137, 2).Value = $docTitle
$ws.Cells.Item(137, 3).Value = $orgNum
$ws.Cells.Item(137, 4).Value = $orgText
$ws.Cells.Item(137, 5).Value = $orgEnt
$ws.Cells.Item(137, 6).Value = "c."
$ws.Cells.Item(137, 7).Value = "In accordance with Reference (g), implement the standards and procedures governing the operation, management, and oversight of clinical laboratory assets assigned to operational forces. Except where operational constraints preclude compliance, the standards governing clinical laboratory assets assigned to operational forces will incorporate the CLIP policy to the maximum extent possible without impeding operational requirements."
$ws.Cells.Item(137, 9).Value = "active"

# --- Row 138 -------------------------------------------------------------
$ws.Cells.Item(138, 1).Value = $filename
$ws.Cells.Item(138, 2).Value = $docTitle
$ws.Cells.Item(138, 3).Value = $orgNum
$ws.Cells.Item(138, 4).Value = $orgText
$ws.Cells.Item(138, 5).Value = $orgEnt
$ws.Cells.Item(138, 6).Value = "d."
$ws.Cells.Item(138, 7).Value = "Recommend changes and revisions to CLIP standards to CCLM. ENCLOSURE 2"
$ws.Cells.Item(138, 9).Value = "active"

# --- Row 139 -------------------------------------------------------------
$ws.Cells.Item(139, 1).Value = $filename
$ws.Cells.Item(139, 2).Value = $docTitle
$ws.Cells.Item(139, 3).Value = $orgNum
$ws.Cells.Item(139, 4).Value = $orgText
$ws.Cells.Item(139, 5).Value = $orgEnt
$ws.Cells.Item(139, 6).Value = "e."
$ws.Cells.Item(139, 7).Value = "Oversee Surgeon General, laboratory commander, and laboratory medical director implementation of the procedures in Enclosure 3 of this instruction."
$ws.Cells.Item(139, 9).Value = "active"
